$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.287671089172363
$ws.Range("B1").Value = 1.481428265571594
$ws.Range("C1").Value = 3.746825933456421
$ws.Range("D1").Value = 3.471381902694702
$ws.Range("E1").Value = 1.008962273597717
